$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Professional summary: neutralize "all Black and Asian-American voters"
#    -> "50M voters" (plain text run, no formatting change)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "affecting all Black and Asian-American voters, developed geospatial ML",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M voters, developed geospatial ML", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Siege Analytics bullet: same neutralization, but "50M" must land in its
#    own bold / colored run (matching the "23%"/"64%" stat runs nearby).
#    Replace "all Black and Asian-American" with "50M" (keeping " voters,
#    developed..." intact), then format just that replaced run.
# ---------------------------------------------------------------------------
$rngBullet = $d.Content
$rngBullet.Find.Execute(
    "all Black and Asian-American",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngBullet.Text = "50M"
$rngBullet.Font.Bold = 1
$rngBullet.Font.Color = 5258796   # RGB(44,62,80) == hex 2C3E50

# ---------------------------------------------------------------------------
# 3) Move the "Software Engineer - Mautinoa Technologies" block so it follows
#    the "Partner - Siege Analytics" section instead of the "Software
#    Engineer - Salsa Labs" section.
# ---------------------------------------------------------------------------
$mautinoaHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Software Engineer - Mautinoa Technologies")) {
        $mautinoaHeading = $p
        break
    }
}

$startPara = $mautinoaHeading
$endPara = $startPara
# The block is the heading paragraph plus the following non-heading
# paragraphs (sub-title + bullet points), i.e. everything up to (but not
# including) the next heading of any level.
while ($true) {
    $nextPara = $endPara.Next()
    if ($nextPara -eq $null) { break }
    if ($nextPara.Style.NameLocal -eq "Heading 3" -or $nextPara.Style.NameLocal -eq "Heading 2") { break }
    $endPara = $nextPara
}

$moveRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$moveRange.Cut()

$myersHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Senior Analyst - Myers Research")) {
        $myersHeading = $p
        break
    }
}
$insertionPoint = $d.Range($myersHeading.Range.Start, $myersHeading.Range.Start)
$insertionPoint.Paste()

# The paste can drop the heading paragraph's style, so make sure the moved
# heading keeps its "Heading 3" style.
$movedHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Software Engineer - Mautinoa Technologies")) {
        $movedHeading = $p
        break
    }
}
$movedHeading.Style = "Heading 3"

# ---------------------------------------------------------------------------
# 4) Key-projects impact line: same neutralization (adds "nationwide")
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "affecting all Black and Asian-American voters, improved electoral prediction accuracy by 22%",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M voters nationwide, improved electoral prediction accuracy by 22%", 2) | Out-Null

Write-Output "done"
